# Unit 2 Research Questions Drill - apply commit "Research Questions - Good or Bad?"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: First "Good" paragraph (after "What can we do to reduce juvenile
# delinquency in the U.S.?") becomes a "Bad:" explanation.
# ---------------------------------------------------------------------------
$pGood = $d.Paragraphs(5)
if ($pGood.Range.Text.TrimEnd([char]13) -ne "Good") {
    throw "Unexpected paragraph 5 content: $($pGood.Range.Text)"
}
$null = $pGood.Range.Find.Execute("Good", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Bad: A better question would be more specific: Is there a correlation between juvenile delinquency and number of parents in the home? ", 2)

# ---------------------------------------------------------------------------
# Change 2: "ATT's " -> "AT" + "&" + "T's " (three separate runs, same
# formatting) inside "A better research question would be: ATT's current
# marketing efforts ..."
# ---------------------------------------------------------------------------
$pAtt = $d.Paragraphs(10)
$attRange = $pAtt.Range.Duplicate
$null = $attRange.Find.Execute("ATT's ", $true)
$attStart = $attRange.Start

# Isolate the whole "ATT's " span into its own run (split it off from the
# preceding "... What impact will " run) by toggling a format and back.
$attRange.Bold = 1
$attRange.Bold = 0

# Insert the ampersand between "AT" and "T's ".
$insertPoint = $d.Range($attStart + 2, $attStart + 2)
$insertPoint.InsertAfter("&")
# Isolate the inserted "&" into its own run too.
$ampRange = $d.Range($attStart + 2, $attStart + 3)
$ampRange.Bold = 1
$ampRange.Bold = 0

# ---------------------------------------------------------------------------
# Change 3: Add <w:lastRenderedPageBreak/> right before "Why did the
# Challenger Shuttle explode?" and delete the following empty ListParagraph.
# ---------------------------------------------------------------------------
$pChallenger = $d.Paragraphs(15)
if ($pChallenger.Range.Text.TrimEnd([char]13) -ne "Why did the Challenger Shuttle explode?") {
    throw "Unexpected paragraph 15 content: $($pChallenger.Range.Text)"
}
$pbPoint = $d.Range($pChallenger.Range.Start, $pChallenger.Range.Start)
$pbXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="373A36"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pbPoint.InsertXML($pbXml)

# The empty ListParagraph-styled paragraph right after "Why did the
# Challenger Shuttle explode?" is removed entirely.
$pEmptyListPara = $d.Paragraphs(16)
if ($pEmptyListPara.Range.Text.TrimEnd([char]13) -ne "") {
    throw "Unexpected paragraph 16 content: $($pEmptyListPara.Range.Text)"
}
$pEmptyListPara.Range.Delete()

# ---------------------------------------------------------------------------
# Change 4: After "... How does Google determine the ranking of their
# indexed web-pages?" add a new paragraph: "Or: " + (moved _GoBack bookmark)
# + "What are the most important factors in Google's algorithm to obtain a
# high search rank?"
# ---------------------------------------------------------------------------
$pGoogleRank = $d.Paragraphs(23)
if ($pGoogleRank.Range.Text -notlike "*How does Google determine the ranking of their indexed web-pages?*") {
    throw "Unexpected paragraph 23 content: $($pGoogleRank.Range.Text)"
}
$pGoogleRank.Range.InsertParagraphAfter()
$pOr = $d.Paragraphs(24)
$pOr.Range.Text = "Or: What are the most important factors in Google" + [char]8217 + "s algorithm to obtain a high search rank?"

# Split "Or: " into its own run.
$orSplit = $pOr.Range.Duplicate
$null = $orSplit.Find.Execute("Or: ", $true)
$orSplit.Bold = 1
$orSplit.Bold = 0

# Move the _GoBack bookmark from the end of the previous paragraph to
# between "Or: " and the new question text.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$orAgain = $pOr.Range.Duplicate
$null = $orAgain.Find.Execute("Or: ", $true)
$bookmarkPoint = $d.Range($orAgain.End, $orAgain.End)
$null = $d.Bookmarks.Add("_GoBack", $bookmarkPoint)

Write-Output "Done"
